# CI list - adding links
# Applies the commit: renames a UML item, clears a stale link, fixes two
# label typos, and wires up hyperlinks for every "Repository" cell that
# was missing one (pointing at the Team-Four GitHub repo / Pivotal Tracker /
# the project website, matching the existing links already on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fixes -------------------------------------------------------
$ws.Range("B6").Value = "Use Case UML"
$ws.Range("F12").Value = "GitHub"
$ws.Range("F19").Value = "Loans and Grants Website"
$ws.Range("F21").Value = "Pivotal Tracker"

# F11 loses its stray "GitHub" text (no link ever existed for this row).
$ws.Range("F11").Value = ""

# --- New hyperlinks -----------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("F6"), "https://github.com/wshahzad/Team-Four/blob/master/Documents/Use%20Case%20UML.docx")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://github.com/wshahzad/Team-Four/blob/master/Documents/UML%20-%20Components%20Diagram.docx")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://github.com/wshahzad/Team-Four/blob/master/Documents/UML%20-%20State%20Diagrams.docx")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://github.com/wshahzad/Team-Four/blob/master/Documents/Use%20cases.docx")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://github.com/wshahzad/Team-Four/blob/master/Documents/Estimation%20Record.docx")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://github.com/wshahzad/Team-Four/blob/master/Documents/Mocks%20ups-Wireframes.docx")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://github.com/wshahzad/Team-Four/blob/master/Documents/Test%20Cases.docx")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.pivotaltracker.com/n/projects/1349388")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://github.com/wshahzad/Team-Four")
$ws.Hyperlinks.Add($ws.Range("F21"), "https://www.pivotaltracker.com/n/projects/1349388")
$ws.Hyperlinks.Add($ws.Range("F22"), "https://github.com/wshahzad/Team-Four/blob/master/Documents/Issues%20from%20Peer%20Reviews%20retained.docx")

# --- Selection moved by the author while reviewing the new links --------
$ws.Range("G12").Select()
